# Scheduled-runner update: refresh computed market-price / profit
# columns (H:N) on the per-class Leve tables. Source data are plain
# numeric literals (no formulas in this workbook), so each corrected
# figure is written straight into its cell via Range.Value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 8930.333000000001
$ws.Range("J64").Value = 10195
$ws.Range("L64").Value = 10195
$ws.Range("N64").Value = -10691
$ws.Range("H67").Value = 8930.333000000001
$ws.Range("J67").Value = 10195
$ws.Range("L67").Value = 10195
$ws.Range("N67").Value = -11911
$ws.Range("H100").Value = 2686.15
$ws.Range("I100").Value = 2901.889
$ws.Range("J100").Value = 2509.6365
$ws.Range("K100").Value = 2901.889
$ws.Range("L100").Value = 2509.6365
$ws.Range("M100").Value = -2360.889
$ws.Range("N100").Value = -3591.6365
$ws.Range("H107").Value = 820.63635
$ws.Range("I107").Value = 718.0625
$ws.Range("J107").Value = 1094.1666
$ws.Range("K107").Value = 718.0625
$ws.Range("L107").Value = 1094.1666
$ws.Range("M107").Value = 1201.9375
$ws.Range("N107").Value = -4934.1666
$ws.Range("H116").Value = 3442.077
$ws.Range("I116").Value = 3326
$ws.Range("J116").Value = 3929.6
$ws.Range("K116").Value = 3326
$ws.Range("L116").Value = 3929.6
$ws.Range("M116").Value = 116
$ws.Range("N116").Value = -10813.6
$ws.Range("H132").Value = 3028.6956
$ws.Range("I132").Value = 3165.1843
$ws.Range("K132").Value = 9495.552899999999
$ws.Range("M132").Value = -6965.552899999999
$ws.Range("H137").Value = 1586.6957
$ws.Range("I137").Value = 1531.421
$ws.Range("J137").Value = 1625.5927
$ws.Range("K137").Value = 4594.263
$ws.Range("L137").Value = 4876.7781
$ws.Range("M137").Value = -2044.263
$ws.Range("N137").Value = -9976.7781
$ws.Range("H138").Value = 6929.4736
$ws.Range("J138").Value = 7053.7144
$ws.Range("L138").Value = 21161.1432
$ws.Range("N138").Value = -31441.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2559.0557
$ws.Range("I2").Value = 2185.5715
$ws.Range("K2").Value = 2185.5715
$ws.Range("M2").Value = -2072.5715
$ws.Range("H25").Value = 1660.75
$ws.Range("I25").Value = 1214.3334
$ws.Range("J25").Value = 3000
$ws.Range("K25").Value = 1214.3334
$ws.Range("L25").Value = 3000
$ws.Range("M25").Value = -812.3334
$ws.Range("N25").Value = -3804
$ws.Range("H32").Value = 6740.6553
$ws.Range("I32").Value = 689.2449
$ws.Range("J32").Value = 39687.223
$ws.Range("K32").Value = 689.2449
$ws.Range("L32").Value = 39687.223
$ws.Range("M32").Value = -402.2449
$ws.Range("N32").Value = -40261.223
$ws.Range("H45").Value = 2691.484
$ws.Range("I45").Value = 1969
$ws.Range("K45").Value = 1969
$ws.Range("M45").Value = -1592
$ws.Range("H74").Value = 2012.8667
$ws.Range("I74").Value = 1131.1111
$ws.Range("K74").Value = 1131.1111
$ws.Range("M74").Value = -257.1111000000001
$ws.Range("H77").Value = 2012.8667
$ws.Range("I77").Value = 1131.1111
$ws.Range("K77").Value = 5655.5555
$ws.Range("M77").Value = -1287.5555
$ws.Range("H110").Value = 2425.0715
$ws.Range("I110").Value = 2380.8462
$ws.Range("K110").Value = 2380.8462
$ws.Range("M110").Value = -335.8462
$ws.Range("H116").Value = 2559.0557
$ws.Range("I116").Value = 2185.5715
$ws.Range("K116").Value = 2185.5715
$ws.Range("M116").Value = 108.4285
$ws.Range("H132").Value = 3306.5642
$ws.Range("I132").Value = 3250.2122
$ws.Range("K132").Value = 9750.6366
$ws.Range("M132").Value = -7220.6366

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2559.0557
$ws.Range("I3").Value = 2185.5715
$ws.Range("K3").Value = 2185.5715
$ws.Range("M3").Value = -2071.5715
$ws.Range("H80").Value = 3074.5
$ws.Range("I80").Value = 2649.3333
$ws.Range("J80").Value = 3499.6667
$ws.Range("K80").Value = 2649.3333
$ws.Range("L80").Value = 3499.6667
$ws.Range("M80").Value = -1651.3333
$ws.Range("N80").Value = -5495.6667
$ws.Range("H83").Value = 3074.5
$ws.Range("I83").Value = 2649.3333
$ws.Range("J83").Value = 3499.6667
$ws.Range("K83").Value = 13246.6665
$ws.Range("L83").Value = 17498.3335
$ws.Range("M83").Value = -8254.666499999999
$ws.Range("N83").Value = -27482.3335
$ws.Range("H94").Value = 628.9
$ws.Range("I94").Value = 587.7778
$ws.Range("K94").Value = 587.7778
$ws.Range("M94").Value = -136.7778
$ws.Range("H99").Value = 4160.2
$ws.Range("I99").Value = 3950.25
$ws.Range("K99").Value = 3950.25
$ws.Range("M99").Value = -2452.25
$ws.Range("H107").Value = 2312.5833
$ws.Range("I107").Value = 1996.1177
$ws.Range("J107").Value = 3081.1428
$ws.Range("K107").Value = 1996.1177
$ws.Range("L107").Value = 3081.1428
$ws.Range("M107").Value = -76.11770000000001
$ws.Range("N107").Value = -6921.1428
$ws.Range("H134").Value = 1834.5454
$ws.Range("I134").Value = 1696.2222
$ws.Range("K134").Value = 5088.6666
$ws.Range("M134").Value = -2553.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5812.778
$ws.Range("I31").Value = 6111.0415
$ws.Range("K31").Value = 6111.0415
$ws.Range("M31").Value = -5816.0415
$ws.Range("H34").Value = 5812.778
$ws.Range("I34").Value = 6111.0415
$ws.Range("K34").Value = 6111.0415
$ws.Range("M34").Value = -5909.0415
$ws.Range("H132").Value = 5111.9414
$ws.Range("I132").Value = 5088.4
$ws.Range("K132").Value = 15265.2
$ws.Range("M132").Value = -12735.2
$ws.Range("H134").Value = 6040.5713
$ws.Range("I134").Value = 4556.9
$ws.Range("K134").Value = 13670.7
$ws.Range("M134").Value = -11135.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2923.6
$ws.Range("I3").Value = 2282.4285
$ws.Range("K3").Value = 6847.2855
$ws.Range("M3").Value = -6735.2855
$ws.Range("H92").Value = 379
$ws.Range("I92").Value = 397.5
$ws.Range("J92").Value = 366.66666
$ws.Range("K92").Value = 1192.5
$ws.Range("L92").Value = 1099.99998
$ws.Range("M92").Value = 55.5
$ws.Range("N92").Value = -3595.99998
$ws.Range("H131").Value = 2137.4
$ws.Range("J131").Value = 2345.6
$ws.Range("L131").Value = 7036.799999999999
$ws.Range("N131").Value = -17116.8
$ws.Range("H138").Value = 8302
$ws.Range("I138").Value = 2453.3333
$ws.Range("K138").Value = 7359.999899999999
$ws.Range("M138").Value = -2219.999899999999
$ws.Range("H141").Value = 8313.25
$ws.Range("I141").Value = 8313.25
$ws.Range("K141").Value = 24939.75
$ws.Range("M141").Value = -19759.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3853.7144
$ws.Range("I122").Value = 3217.8572
$ws.Range("J122").Value = 5125.4287
$ws.Range("K122").Value = 9653.571599999999
$ws.Range("L122").Value = 15376.2861
$ws.Range("M122").Value = -7203.571599999999
$ws.Range("N122").Value = -20276.2861
$ws.Range("H132").Value = 3370.375
$ws.Range("I132").Value = 3586.158
$ws.Range("K132").Value = 10758.474
$ws.Range("M132").Value = -8228.474

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19153.5
$ws.Range("I7").Value = 13913.385
$ws.Range("J7").Value = 32777.8
$ws.Range("K7").Value = 13913.385
$ws.Range("L7").Value = 32777.8
$ws.Range("M7").Value = -13801.385
$ws.Range("N7").Value = -33001.8
$ws.Range("H40").Value = 6369.476
$ws.Range("I40").Value = 3466.9167
$ws.Range("K40").Value = 3466.9167
$ws.Range("M40").Value = -3330.9167
$ws.Range("H46").Value = 4085.875
$ws.Range("I46").Value = 1250.5
$ws.Range("K46").Value = 1250.5
$ws.Range("M46").Value = -1062.5
$ws.Range("H93").Value = 1532.25
$ws.Range("I93").Value = 1497.6
$ws.Range("J93").Value = 1590
$ws.Range("K93").Value = 1497.6
$ws.Range("L93").Value = 1590
$ws.Range("M93").Value = -249.5999999999999
$ws.Range("N93").Value = -4086
$ws.Range("H126").Value = 19153.5
$ws.Range("I126").Value = 13913.385
$ws.Range("J126").Value = 32777.8
$ws.Range("K126").Value = 41740.155
$ws.Range("L126").Value = 98333.40000000001
$ws.Range("M126").Value = -39270.155
$ws.Range("N126").Value = -103273.4
$ws.Range("H132").Value = 23286.973
$ws.Range("I132").Value = 35575.863
$ws.Range("K132").Value = 106727.589
$ws.Range("M132").Value = -104197.589
$ws.Range("H136").Value = 6329.5
$ws.Range("I136").Value = 5369.375
$ws.Range("K136").Value = 16108.125
$ws.Range("M136").Value = -13558.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2671.0386
$ws.Range("I107").Value = 3413.4546
$ws.Range("K107").Value = 10240.3638
$ws.Range("M107").Value = -8320.363799999999
$ws.Range("H126").Value = 3312
$ws.Range("I126").Value = 2614
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 7842
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -5372
$ws.Range("N126").Value = -27440
$ws.Range("H132").Value = 2321.6458
$ws.Range("I132").Value = 2069.2632
$ws.Range("J132").Value = 3280.7
$ws.Range("K132").Value = 6207.7896
$ws.Range("L132").Value = 9842.099999999999
$ws.Range("M132").Value = -3677.7896
$ws.Range("N132").Value = -14902.1
